$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so that values like
# "24.918.21", "1.004", "0.00001139" etc. are preserved exactly as strings
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 30/31 swap: ImmutableX <-> WrappedliquidstakedEther2.0 (with updated values)
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.903.12"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.254"
$ws.Range("E31").Value = "  +28.16%  "

# Row 44/45 swap: EnergySwap <-> Frax (with updated values)
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "0.9979"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.07"
$ws.Range("E45").Value = "  +6.11%  "

# Remaining D/E value updates for all other rows
$ws.Range("D2").Value = "24.918.21"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.712.27"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "313.48"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "49.72"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "0.3473"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "1.219"
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("D11").Value = "0.07604"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "21.44"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").Value = "6.365"
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("D15").Value = "7.093"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("D16").Value = "1.713.44"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "0.00001139"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("D18").Value = "0.06753"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "0.9976"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "85.30"
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("D21").Value = "17.44"
$ws.Range("E21").Value = "  +5.97%  "
$ws.Range("D22").Value = "6.431"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").Value = "13.30"
$ws.Range("E23").Value = "  +10.31%  "
$ws.Range("D24").Value = "24.931.22"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "2.457"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "2.814"
$ws.Range("E26").Value = "  +4.74%  "
$ws.Range("D27").Value = "20.59"
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("D28").Value = "151.43"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "133.41"
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D32").Value = "6.941"
$ws.Range("E32").Value = "  +9.05%  "
$ws.Range("D33").Value = "4.249"
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("D34").Value = "14.01"
$ws.Range("E34").Value = "  +11.04%  "
$ws.Range("D35").Value = "1.794"
$ws.Range("E35").Value = "  +5.57%  "
$ws.Range("D36").Value = "0.08879"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("D37").Value = "5.684"
$ws.Range("E37").Value = "  +5.63%  "
$ws.Range("D38").Value = "9.413"
$ws.Range("E38").Value = "  +5.31%  "
$ws.Range("D39").Value = "0.06710"
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("D40").Value = "0.02429"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("D41").Value = "0.2253"
$ws.Range("E41").Value = "  +6.50%  "
$ws.Range("D42").Value = "1.286"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").Value = "0.6501"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D46").Value = "0.6203"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "3.847"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").Value = "2.151"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("D49").Value = "130.83"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("D50").Value = "0.07328"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "80.62"
$ws.Range("E51").Value = "  +6.00%  "

Write-Output "Applied cryptos update"
